$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "DE Arm 25 2.5mm Parts" with "DE Arm 25 2mm Parts" in A5.
$ws.Range("A5").Value = "DE Arm 25 2mm Parts"

# Match the bold style used by the other part-name cells (e.g. A6) by
# copying its formatting onto A5.
$ws.Range("A6").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to A5.
[void]$ws.Range("A5").Select()
